# Generate Report for Handoff
# Replace the old UUID-based file identifiers with the new UUID, and
# refresh the associated handoff timestamps across the three sheets
# (Overview, zh-cn, de-de).

$oldGuid = "1130f50d-6537-492d-ad09-1677807fb620"
$newGuid = "d120e4fd-5ce2-49e9-aba3-09a2b49ec8aa"

$oldHash = "04b4020ae91b5ebb758f4c9e190bdbc1bd19c5b0"
$newHash = "db2273b2c4a71febf59e71415699cafeaac3f6ce"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: File Name
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: Path And Name (hyperlink) - keep the existing link target, only the
# displayed text changes.
$oldHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/122c72732cb908adc2f6ef8e5f03e3bb50632406/e2e/$oldGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $oldHyperlinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")

# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-03 13:02:54"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2: Source File Name (hyperlink) - keep existing link target, only text.
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $oldHyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

# G2: Latest Handoff File
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-03 13:02:50"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2: Source File Name (hyperlink) - keep existing link target, only text.
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $oldHyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

# G2: Latest Handoff File
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# H2: Latest Handback DateTime - shares the same underlying text as
# Overview!G2, so update it the same way.
$wsDeDe.Range("H2").Value = "2016-09-03 13:02:54"
